# Applies the commit's changes to the "EVALUABLE T6 Enunciado" document:
#   1. Collapse the two leading empty paragraphs (one holding a lone
#      <w:br/> run) into a single empty paragraph and drop Word's
#      "_GoBack" last-edit bookmark there.
#   2. Remove the old "_GoBack" bookmark that used to sit in the middle
#      of "...los perm|isos..." and re-merge that sentence into one
#      continuous run of text.
#   3. Drop the stale <w:lastRenderedPageBreak/> marker in front of the
#      final "Jorge" signature line.

$d = $word.ActiveDocument

# --- 1. Merge the first two (empty) paragraphs into one -------------------
$firstPara = $d.Paragraphs(1).Range
$mergeRange = $d.Range($firstPara.Start, $firstPara.End)
$mergeRange.Delete()

# Drop Word's "_GoBack" bookmark at the start of the now-merged paragraph.
# (Adding it here also removes the stale copy left over from the old
# edit location further down in the document, matching real Word's
# single-instance "_GoBack" behaviour.)
$survivingPara = $d.Paragraphs(1).Range
$goBackRange = $d.Range($survivingPara.Start, $survivingPara.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 2. Re-join "...los perm" + "isos..." into one continuous run ---------
$search = "de las carpetas, los permisos como hemos comentado funcionan"
$d.Content.Find.Execute($search, $true, $false, $false, $false, $false, `
    $true, 1, $false, $search, 2) | Out-Null

# --- 3. Remove the stray lastRenderedPageBreak before "Jorge" -------------
$d.Content.Find.Execute("Jorge", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Jorge", 2) | Out-Null
